# Update the "Price" (column D) values for the symbol list refresh,
# matching the GitHub Actions commit on Thu Dec 15 09:43:19 UTC 2022.
# Cells are stored as text (not numbers), so we use a leading apostrophe
# to force a literal/text assignment and keep the exact textual value
# (e.g. "265.90" rather than being coerced to the number 265.9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "265.90"
    "D3"  = "22.63"
    "D4"  = "6.199"
    "D7"  = "6.713"
    "D8"  = "1.356"
    "D9"  = "0.8261"
    "D11" = "0.1590"
    "D12" = "0.08212"
    "D13" = "0.03399"
    "D14" = "0.03155"
    "D15" = "0.09237"
    "D16" = "3.895"
    "D17" = "0.001693"
    "D18" = "0.04813"
    "D19" = "0.006276"
    "D20" = "0.006279"
    "D21" = "0.001098"
    "D22" = "0.0001420"
    "D23" = "3.717"
    "D24" = "2.231"
    "D25" = "0.3378"
    "D40" = "0.04613"
    "D41" = "0.006989"
    "D44" = "0.01103"
    "D47" = "0.7700"
    "D48" = "0.2067"
    "D49" = "0.00002101"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
